# Applies the "Added a few more slots" edit:
#  1. Replace the text of the trailing italic paragraph (previously the
#     "Read our unbiased review..." meta description) with the new
#     image-generation prompt, while keeping its italic formatting intact.
#  2. Remove the duplicate bold "Play Big Max 77 Free..." paragraph that
#     used to sit near the end of the document (its content now lives in
#     the new Meta description paragraph at the top).
#  3. Insert a new "Meta description" paragraph right after the title
#     (Heading1) paragraph. "Meta description" is bold, the rest of the
#     sentence is regular text.
#
# NOTE: step 1 is done first, while the "Read our unbiased review..."
# text still only occurs once in the document - this keeps the
# Find/Replace unambiguous (step 3 later inserts that same sentence
# again, near the top, as part of the Meta description paragraph).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: swap the text of the final (italic) paragraph for the new
#         image-generation brief. Find/Replace turns straight quotes
#         into curly ones, so the quotes are inserted via a plain-text
#         Range.Text assignment afterwards, which keeps them straight.
# ---------------------------------------------------------------------
$oldBlurb = "Read our unbiased review of Big Max 77, the classic fruit-themed slot game with high volatility and an RTP of 97.01%. Play for free and learn more!"
$newBlurb = 'Create a feature image for Big Max 77 to use on social media and marketing materials. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a giant fruit, such as a watermelon or pineapple, and there should be a slot machine in the background with the name ~Big Max 77~ displayed prominently. The overall style should be fun and engaging, with bright colors and playful designs to attract potential players.'

$blurbRange = $d.Content
$blurbRange.Find.Execute($oldBlurb, $true, $false, $false, $false, $false, $true, 1, $false, $newBlurb, 2) | Out-Null

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$quoteSearch = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$quoteSearch.Find.Execute("~Big Max 77~", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$quoteChar = [char]34
$quoteSearch.Text = $quoteChar + "Big Max 77" + $quoteChar

# ---------------------------------------------------------------------
# Step 2: delete the old bold title paragraph further down the document
#         (second occurrence of the title text).
# ---------------------------------------------------------------------
$titleText = "Play Big Max 77 Free - Review & Demo | RTP 97.01%"

$firstHit = $d.Content
$firstHit.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$afterFirst = $d.Range($firstHit.End, $d.Content.End)
$afterFirst.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$dupPara = $afterFirst.Paragraphs.Item(1)
$dupRange = $d.Range($dupPara.Range.Start, $dupPara.Range.End + 1)
$dupRange.Delete()

# ---------------------------------------------------------------------
# Step 3: insert the new "Meta description" paragraph after the title.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"
$metaPara.Range.Text = "Meta description: Read our unbiased review of Big Max 77, the classic fruit-themed slot game with high volatility and an RTP of 97.01%. Play for free and learn more!"

$labelLength = "Meta description".Length
$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + $labelLength)
$boldRange.Font.Bold = 1
